$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 102
$ws.Range("I5").Value = 128
$ws.Range("K5").Value = 128
$ws.Range("M5").Value = -13

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2703.1428
$ws.Range("I34").Value = 2703.1428
$ws.Range("K34").Value = 2703.1428
$ws.Range("M34").Value = -2500.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 2703.1428
$ws.Range("I36").Value = 2703.1428
$ws.Range("K36").Value = 2703.1428
$ws.Range("M36").Value = -1988.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1573.7455
$ws.Range("I112").Value = 1259.5
$ws.Range("J112").Value = 1585.6038
$ws.Range("K112").Value = 3778.5
$ws.Range("L112").Value = 4756.811400000001
$ws.Range("M112").Value = -2670.5
$ws.Range("N112").Value = -6972.811400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 12353.777
$ws.Range("I125").Value = 1191.3334
$ws.Range("K125").Value = 10722.0006
$ws.Range("M125").Value = -8262.000599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6866.6665
$ws.Range("I141").Value = 7128.5713
$ws.Range("K141").Value = 21385.7139
$ws.Range("M141").Value = -16205.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11333.286
$ws.Range("I32").Value = 8901.138999999999
$ws.Range("J32").Value = 18068.46
$ws.Range("K32").Value = 8901.138999999999
$ws.Range("L32").Value = 18068.46
$ws.Range("M32").Value = -8614.138999999999
$ws.Range("N32").Value = -18642.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6848.2856
$ws.Range("I45").Value = 4184.8
$ws.Range("K45").Value = 4184.8
$ws.Range("M45").Value = -3807.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2401.9707
$ws.Range("I61").Value = 1747.1724
$ws.Range("J61").Value = 6199.8
$ws.Range("K61").Value = 1747.1724
$ws.Range("L61").Value = 6199.8
$ws.Range("M61").Value = -1535.1724
$ws.Range("N61").Value = -6623.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2828.6155
$ws.Range("I74").Value = 2647.6667
$ws.Range("K74").Value = 2647.6667
$ws.Range("M74").Value = -1773.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2828.6155
$ws.Range("I77").Value = 2647.6667
$ws.Range("K77").Value = 13238.3335
$ws.Range("M77").Value = -8870.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2824.6924
$ws.Range("I102").Value = 2564.5454
$ws.Range("K102").Value = 2564.5454
$ws.Range("M102").Value = -942.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4347.8423
$ws.Range("I132").Value = 3532.9167
$ws.Range("K132").Value = 10598.7501
$ws.Range("M132").Value = -8068.750100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2401.9707
$ws.Range("I136").Value = 1747.1724
$ws.Range("J136").Value = 6199.8
$ws.Range("K136").Value = 5241.5172
$ws.Range("L136").Value = 18599.4
$ws.Range("M136").Value = -2691.5172
$ws.Range("N136").Value = -23699.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1075.7
$ws.Range("J64").Value = 1034.7778
$ws.Range("L64").Value = 1034.7778
$ws.Range("N64").Value = -1484.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1075.7
$ws.Range("J67").Value = 1034.7778
$ws.Range("L67").Value = 1034.7778
$ws.Range("N67").Value = -2594.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2873.6
$ws.Range("J86").Value = 5573.364
$ws.Range("L86").Value = 5573.364
$ws.Range("N86").Value = -7819.364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2873.6
$ws.Range("J89").Value = 5573.364
$ws.Range("L89").Value = 27866.82
$ws.Range("N89").Value = -39098.82

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4689.9287
$ws.Range("I94").Value = 4054.9167
$ws.Range("J94").Value = 8500
$ws.Range("K94").Value = 4054.9167
$ws.Range("L94").Value = 8500
$ws.Range("M94").Value = -3603.9167
$ws.Range("N94").Value = -9402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 18125.115
$ws.Range("I99").Value = 20710.592
$ws.Range("J99").Value = 3905
$ws.Range("K99").Value = 20710.592
$ws.Range("L99").Value = 3905
$ws.Range("M99").Value = -19212.592
$ws.Range("N99").Value = -6901

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3070.7917
$ws.Range("I105").Value = 3063.5908
$ws.Range("K105").Value = 3063.5908
$ws.Range("M105").Value = -1316.5908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 27495
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 49990
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 49990
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -50846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 45999
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 4642.92
$ws.Range("I107").Value = 702.38464
$ws.Range("J107").Value = 8911.833000000001
$ws.Range("K107").Value = 702.38464
$ws.Range("L107").Value = 8911.833000000001
$ws.Range("M107").Value = 1217.61536
$ws.Range("N107").Value = -12751.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3527.7307
$ws.Range("I132").Value = 3498.6316
$ws.Range("J132").Value = 3606.7144
$ws.Range("K132").Value = 10495.8948
$ws.Range("L132").Value = 10820.1432
$ws.Range("M132").Value = -7965.8948
$ws.Range("N132").Value = -15880.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3943.4583
$ws.Range("I134").Value = 2733.8276
$ws.Range("K134").Value = 8201.4828
$ws.Range("M134").Value = -5666.4828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.8
$ws.Range("I2").Value = 4.6666665
$ws.Range("J2").Value = 77.5
$ws.Range("K2").Value = 27.999999
$ws.Range("L2").Value = 465
$ws.Range("M2").Value = 85.000001
$ws.Range("N2").Value = -691

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1589
$ws.Range("I48").Value = 618.6667
$ws.Range("J48").Value = 4500
$ws.Range("K48").Value = 1856.0001
$ws.Range("L48").Value = 13500
$ws.Range("M48").Value = -1606.0001
$ws.Range("N48").Value = -14000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2750
$ws.Range("J59").Value = 4500
$ws.Range("L59").Value = 13500
$ws.Range("N59").Value = -14580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 69690
$ws.Range("J51").Value = 69690
$ws.Range("L51").Value = 69690
$ws.Range("N51").Value = -70708

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4701.923
$ws.Range("I126").Value = 3337.2856
$ws.Range("K126").Value = 10011.8568
$ws.Range("M126").Value = -7541.856800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4265.7856
$ws.Range("I132").Value = 2978.7778
$ws.Range("K132").Value = 8936.3334
$ws.Range("M132").Value = -6406.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7335.65
$ws.Range("I7").Value = 10915.857
$ws.Range("K7").Value = 10915.857
$ws.Range("M7").Value = -10803.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2177.077
$ws.Range("I16").Value = 2631.5
$ws.Range("J16").Value = 662.3333
$ws.Range("K16").Value = 2631.5
$ws.Range("L16").Value = 662.3333
$ws.Range("M16").Value = -2461.5
$ws.Range("N16").Value = -1002.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 73571.664
$ws.Range("J36").Value = 73571.664
$ws.Range("L36").Value = 73571.664
$ws.Range("N36").Value = -74695.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 108010
$ws.Range("I100").Value = 373036.34
$ws.Range("J100").Value = 8625.125
$ws.Range("K100").Value = 373036.34
$ws.Range("L100").Value = 8625.125
$ws.Range("M100").Value = -372495.34
$ws.Range("N100").Value = -9707.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 59666.582
$ws.Range("J123").Value = 59666.582
$ws.Range("L123").Value = 59666.582
$ws.Range("N123").Value = -69466.58199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 499999
$ws.Range("J124").Value = 499999
$ws.Range("L124").Value = 499999
$ws.Range("N124").Value = -509819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 69999.89999999999
$ws.Range("J125").Value = 69999.89999999999
$ws.Range("L125").Value = 69999.89999999999
$ws.Range("N125").Value = -79839.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7335.65
$ws.Range("I126").Value = 10915.857
$ws.Range("K126").Value = 32747.571
$ws.Range("M126").Value = -30277.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4038.9092
$ws.Range("I132").Value = 3038
$ws.Range("K132").Value = 9114
$ws.Range("M132").Value = -6584

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4254.9
$ws.Range("I136").Value = 2591.4443
$ws.Range("K136").Value = 7774.3329
$ws.Range("M136").Value = -5224.3329
